# Scheduled runner update: refresh market-board derived columns (H-N) across leve profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 44
$ws.Range("H44").Value = 9666.666999999999
$ws.Range("J44").Value = 9666.666999999999
$ws.Range("L44").Value = 9666.666999999999
$ws.Range("N44").Value = -10590.667

# Row 61
$ws.Range("H61").Value = 73.8
$ws.Range("I61").Value = 73.8
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 221.4
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -49.39999999999998
$ws.Range("N61").ClearContents()

# Row 113
$ws.Range("H113").Value = 4358.4116
$ws.Range("I113").Value = 3215.5
$ws.Range("J113").Value = 4981.8184
$ws.Range("K113").Value = 3215.5
$ws.Range("L113").Value = 4981.8184
$ws.Range("M113").Value = 38.5
$ws.Range("N113").Value = -11489.8184

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 752
$ws.Range("I2").Value = 714.1539
$ws.Range("J2").Value = 875
$ws.Range("K2").Value = 714.1539
$ws.Range("L2").Value = 875
$ws.Range("M2").Value = -601.1539
$ws.Range("N2").Value = -1101

# Row 45
$ws.Range("H45").Value = 1190
$ws.Range("I45").Value = 1000
$ws.Range("K45").Value = 1000
$ws.Range("M45").Value = -623

# Row 54
$ws.Range("H54").Value = 9349.75
$ws.Range("J54").Value = 9349.75
$ws.Range("L54").Value = 9349.75
$ws.Range("N54").Value = -10887.75

# Row 116
$ws.Range("H116").Value = 752
$ws.Range("I116").Value = 714.1539
$ws.Range("J116").Value = 875
$ws.Range("K116").Value = 714.1539
$ws.Range("L116").Value = 875
$ws.Range("M116").Value = 1579.8461
$ws.Range("N116").Value = -5463

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 752
$ws.Range("I3").Value = 714.1539
$ws.Range("J3").Value = 875
$ws.Range("K3").Value = 714.1539
$ws.Range("L3").Value = 875
$ws.Range("M3").Value = -600.1539
$ws.Range("N3").Value = -1103

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

# Row 59
$ws.Range("H59").Value = 47743.2
$ws.Range("I59").Value = 40000
$ws.Range("J59").Value = 48603.555
$ws.Range("K59").Value = 40000
$ws.Range("L59").Value = 48603.555
$ws.Range("M59").Value = -38855
$ws.Range("N59").Value = -50893.555

# Row 86
$ws.Range("H86").Value = 5799.625
$ws.Range("I86").Value = 3188.2222
$ws.Range("J86").Value = 9157.143
$ws.Range("K86").Value = 3188.2222
$ws.Range("L86").Value = 9157.143
$ws.Range("M86").Value = -2065.2222
$ws.Range("N86").Value = -11403.143

# Row 89
$ws.Range("H89").Value = 5799.625
$ws.Range("I89").Value = 3188.2222
$ws.Range("J89").Value = 9157.143
$ws.Range("K89").Value = 15941.111
$ws.Range("L89").Value = 45785.715
$ws.Range("M89").Value = -10325.111
$ws.Range("N89").Value = -57017.715

# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

# Row 122
$ws.Range("H122").Value = 1145
$ws.Range("I122").Value = 966.6667
$ws.Range("J122").Value = 1234.1666
$ws.Range("K122").Value = 2900.0001
$ws.Range("L122").Value = 3702.4998
$ws.Range("M122").Value = -450.0001000000002
$ws.Range("N122").Value = -8602.4998

$ws = $wb.Worksheets.Item("CUL")
# Row 62
$ws.Range("H62").Value = 2795.4285
$ws.Range("J62").Value = 3095
$ws.Range("L62").Value = 9285
$ws.Range("N62").Value = -10657

# Row 63
$ws.Range("H63").Value = 3733.3333
$ws.Range("I63").Value = 1000
$ws.Range("J63").Value = 4075
$ws.Range("K63").Value = 3000
$ws.Range("L63").Value = 12225
$ws.Range("M63").Value = -2251
$ws.Range("N63").Value = -13723

# Row 65
$ws.Range("H65").Value = 2795.4285
$ws.Range("J65").Value = 3095
$ws.Range("L65").Value = 27855
$ws.Range("N65").Value = -34719

# Row 66
$ws.Range("H66").Value = 3733.3333
$ws.Range("I66").Value = 1000
$ws.Range("J66").Value = 4075
$ws.Range("K66").Value = 9000
$ws.Range("L66").Value = 36675
$ws.Range("M66").Value = -5256
$ws.Range("N66").Value = -44163

# Row 131
$ws.Range("H131").Value = 725.64
$ws.Range("I131").Value = 307.2
$ws.Range("J131").Value = 904.97144
$ws.Range("K131").Value = 921.5999999999999
$ws.Range("L131").Value = 2714.91432
$ws.Range("M131").Value = 4118.4
$ws.Range("N131").Value = -12794.91432

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 6684.5
$ws.Range("I113").Value = 8595.538
$ws.Range("J113").Value = 1715.8
$ws.Range("K113").Value = 8595.538
$ws.Range("L113").Value = 1715.8
$ws.Range("M113").Value = -6425.538
$ws.Range("N113").Value = -6055.8

# Row 122
$ws.Range("H122").Value = 1340.909
$ws.Range("I122").Value = 1100
$ws.Range("J122").Value = 1394.4445
$ws.Range("K122").Value = 3300
$ws.Range("L122").Value = 4183.333500000001
$ws.Range("M122").Value = -850
$ws.Range("N122").Value = -9083.333500000001

$ws = $wb.Worksheets.Item("LTW")
# Row 81
$ws.Range("H81").Value = 29998.25
$ws.Range("J81").Value = 29998.25
$ws.Range("L81").Value = 29998.25
$ws.Range("N81").Value = -31994.25

# Row 84
$ws.Range("H84").Value = 29998.25
$ws.Range("J84").Value = 29998.25
$ws.Range("L84").Value = 89994.75
$ws.Range("N84").Value = -99978.75

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 14500
$ws.Range("J54").Value = 28000
$ws.Range("L54").Value = 28000
$ws.Range("N54").Value = -29040

# Row 81
$ws.Range("H81").Value = 2412.4375
$ws.Range("I81").Value = 1042.7142
$ws.Range("J81").Value = 3477.7778
$ws.Range("K81").Value = 2085.4284
$ws.Range("L81").Value = 6955.5556
$ws.Range("M81").Value = -1024.4284
$ws.Range("N81").Value = -9077.5556

# Row 84
$ws.Range("H84").Value = 2412.4375
$ws.Range("I84").Value = 1042.7142
$ws.Range("J84").Value = 3477.7778
$ws.Range("K84").Value = 10427.142
$ws.Range("L84").Value = 34777.778
$ws.Range("M84").Value = -5123.142
$ws.Range("N84").Value = -45385.778

# Row 95
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492

# Row 107
$ws.Range("H107").Value = 381.625
$ws.Range("I107").Value = 437.5
$ws.Range("J107").Value = 325.75
$ws.Range("K107").Value = 1312.5
$ws.Range("L107").Value = 977.25
$ws.Range("M107").Value = 607.5
$ws.Range("N107").Value = -4817.25

Write-Host "Updated market price columns on ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."
